# Swap the full contents of row 2 <-> row 4 and row 5 <-> row 6 on the
# "Artfynd" worksheet. Only the cells that actually carry differing data
# between each pair are touched, so unrelated (identical-in-both-rows)
# cells such as the Startdatum/Slutdatum text are left completely alone
# (this avoids Excel's automatic text->date coercion on round trip).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($ws, $colLetter, $row1, $row2) {
    $c1 = $ws.Range("$colLetter$row1")
    $c2 = $ws.Range("$colLetter$row2")
    $v1 = $c1.Value2
    $v2 = $c2.Value2
    $c1.Value2 = $v2
    $c2.Value2 = $v1
}

# --- Row 2 <-> Row 4 : Id / Taxonsorteringsordning / TaxonId / Artnamn /
#     Vetenskapligt namn / Auktor / Ost / Nord all swap; everything else
#     (C,D,I,P,S,T,U,V,W,Y,AA,AD,AE,AG,AT,AW,AX,AY) is identical between
#     the two rows already.
foreach ($col in @("A","B","E","F","G","H","Q","R")) {
    Swap-Cell $ws $col 2 4
}

# --- Row 5 <-> Row 6 : same set of columns swap values ---
foreach ($col in @("A","B","E","F","G","H","Q","R")) {
    Swap-Cell $ws $col 5 6
}

# Row 5 loses its (empty) K/L/M/N placeholder cells - they move to row 6
# (row 6 had none before, and the cells carry no data either way, so the
# net effect of the swap is simply to clear them from row 5).
$ws.Range("K5:N5").ClearContents()

# Row 5 loses its "Publik kommentar" (AC) note - it moves to row 6.
$ws.Range("AC5").ClearContents()
$ws.Range("AC6").Value = "ringhack äldre"
